# Implement calculating vacancy percentage by passing in the total unit
# counts for each building -- refreshed vacancy counts + updated
# availability/pricing snapshot (6/2 -> 6/10) across the competitor sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nl = [char]10

# Row 2 - WILCO Apartments
$ws.Range("B2").Value = "Studio: 3${nl}1 bed: 4${nl}2 bed: 2${nl}Percent: %"
$ws.Range("D2").Value = "6/10: `$2081-2103${nl}(`$3.46-4.91)"
$ws.Range("E2").Value = "6/10: `$2785-3005${nl}(`$3.36-4.19)"
$ws.Range("F2").Value = "6/10: `$3290-3457${nl}(`$3.05-3.37)"

# Row 3 - The Chadwick
$ws.Range("B3").Value = "Studio: 10${nl}1 bed: 6${nl}2 bed: 13${nl}Percent: %"
$ws.Range("D3").Value = "6/10: `$1485-1679${nl}(`$3.55-4.26)"
$ws.Range("E3").Value = "6/10: `$1809-2009${nl}(`$2.59-2.81)"
$ws.Range("F3").Value = "6/10: `$2492-2667${nl}(`$2.27-2.8)"

# Row 4 - The 900
$ws.Range("B4").Value = "Studio: 5${nl}1 bed: 14${nl}2 bed: 3${nl}Percent: %"
$ws.Range("D4").Value = "6/10: `$2200-2550${nl}(`$4-4.48)"
$ws.Range("E4").Value = "6/10: `$2392-3433${nl}(`$3.48-4.62)"
$ws.Range("F4").Value = "6/10: `$3233-3934${nl}(`$3.85-4.16)"

# Row 5 - Fedora x Trilby
$ws.Range("B5").Value = "Studio: 6${nl}1 bed: 3${nl}2 bed: 11${nl}Percent: %"
$ws.Range("D5").Value = "6/10: `$1956-2195${nl}(`$3.81-4.46)"
$ws.Range("E5").Value = "6/10: `$2483-2735${nl}(`$3.51-3.6)"
$ws.Range("F5").Value = "6/10: `$2858-3532${nl}(`$2.76-3.19)"

# Row 6 - Qwil Apartments
$ws.Range("B6").Value = "2 bed: 4${nl}Percent: %"
$ws.Range("F6").Value = "6/10: `$3895-4895${nl}(`$3.53-4.1)"

# Row 7 - The Rhys
$ws.Range("E7").Value = "6/10: `$2425${nl}(`$3.77)"

# Row 8 - Berkshire K2LA
$ws.Range("B8").Value = "Studio: 7${nl}1 bed: 1${nl}Percent: %"
$ws.Range("D8").Value = "6/10: `$1959-2369${nl}(`$4.14-4.7)"
$ws.Range("E8").Value = "6/10: `$2774${nl}(`$4.13)"

# Row 9 - 4749 Elmwood Ave
$ws.Range("B9").Value = "1 bed: 1${nl}2 bed: 2${nl}3 bed: 5${nl}4 bed: 4${nl}Percent: %"
$ws.Range("E9").Value = "6/10: `$2599${nl}(`$3.32)"
$ws.Range("F9").Value = "6/10: `$2599-2949${nl}(`$3.29-5.16)"
$ws.Range("G9").Value = "6/10: `$1198-3750${nl}(`$3.56)"

# Row 10 - Rise Koreatown
$ws.Range("B10").Value = "Studio: 8${nl}1 bed: 13${nl}Percent: %"
$ws.Range("D10").Value = "6/10: `$2405-2530${nl}(`$5.26-5.54)"
$ws.Range("E10").Value = "6/10: `$3250-4550${nl}(`$3.99-6.31)"

# Row 11 - 2783 Francis Ave
$ws.Range("D11").Value = "6/10: `$2125${nl}(`$3.42)"
$ws.Range("E11").Value = "6/10: `$2700-2745${nl}(`$2.74-2.78)"
$ws.Range("F11").Value = "6/10: `$2790-3785${nl}(`$2.76-3.17)"

# Row 12 - Hollywood Flats
$ws.Range("D12").Value = "6/10: `$2245${nl}(`$3.74)"
$ws.Range("E12").Value = "6/10: `$2625-2850${nl}(`$3-3.39)"
$ws.Range("F12").Value = "6/10: `$3450-3675${nl}(`$2.64-3.3)"

# Row 13 - Hallasan
$ws.Range("B13").Value = "Studio: 8${nl}1 bed: 22${nl}2 bed: 32${nl}3 bed: 13${nl}Percent: %"
$ws.Range("D13").Value = "6/10: `$2395-2830${nl}(`$4.57-5.59)"
$ws.Range("E13").Value = "6/10: `$3620-5759${nl}(`$4.11-5.01)"
$ws.Range("F13").Value = "6/10: `$4398-8044${nl}(`$3.73-5.27)"
$ws.Range("G13").Value = "6/10: `$5276-8959${nl}(`$3.67-5)"

# Row 14 - The BORA 3170
$ws.Range("C14").Value = "*Up to 12 Weeks of Free Rent! *1. 13 month lease = 8 weeks free *2. 20 month lease = 12 weeks free *3. `$700 credit towards the rent for look and lease (within 24 hours of the tour) *4. `$50 discount on parking fee each month for the entire lease period (Original fee `$100/month) **All concessions are included in the rent amount."
$ws.Range("D14").Value = "6/10: `$1910-2150${nl}(`$4.47-4.68)"
$ws.Range("E14").Value = "6/10: `$2444-3152${nl}(`$3.89-4.41)"
$ws.Range("F14").Value = "6/10: `$3660${nl}(`$3.74)"

# Row 15 - Miles at Harvard
$ws.Range("B15").Value = "Studio: 13${nl}1 bed: 2${nl}Percent: %"
$ws.Range("D15").Value = "6/10: `$1407-2750${nl}(`$7.07-10.94)"
$ws.Range("E15").Value = "6/10: `$2009-2292"
